$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update wm.sd (E) and wm.mean (F) values for rows 2, 10, 13 with new
# recalculated coefficient statistics (absolute values now used for ES).
$ws.Range("E2").Value = 0.240929464501854
$ws.Range("F2").Value = 14.9050332812142

$ws.Range("E10").Value = 0.240929464501854
$ws.Range("F10").Value = 14.9050332812142

$ws.Range("E13").Value = 0.240929464501854
$ws.Range("F13").Value = 14.9050332812142
